# no-op
$d = $word.ActiveDocument
